$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "36.896.04"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "1.988.63"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue $ws.Range("D5") "241.00"
$ws.Range("E5").Value = "  -4.33%  "
Set-TextValue $ws.Range("D6") "0.603"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  +0.17%  "
Set-TextValue $ws.Range("D8") "54.11"
$ws.Range("E8").Value = "  -4.76%  "
Set-TextValue $ws.Range("D9") "0.372"
$ws.Range("E9").Value = "  -3.41%  "
Set-TextValue $ws.Range("D10") "57.23"
$ws.Range("E10").Value = "  +0.12%  "
Set-TextValue $ws.Range("D11") "0.0750"
$ws.Range("E11").Value = "  -4.60%  "
Set-TextValue $ws.Range("D12") "0.0977"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.284.18"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "14.09"
$ws.Range("E14").Value = "  -3.85%  "
Set-TextValue $ws.Range("D15") "20.80"
$ws.Range("E15").Value = "  -1.54%  "
Set-TextValue $ws.Range("D16") "0.753"
$ws.Range("E16").Value = "  -7.44%  "
Set-TextValue $ws.Range("D17") "5.04"
$ws.Range("E17").Value = "  -5.59%  "
$ws.Range("D18").Value = "2.001.17"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "36.849.84"
$ws.Range("E19").Value = "  -1.45%  "
Set-TextValue $ws.Range("D20") "68.31"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "0.0₃0806"
$ws.Range("E21").Value = "  -4.60%  "
Set-TextValue $ws.Range("D22") "5.08"
$ws.Range("E22").Value = "  -1.62%  "
Set-TextValue $ws.Range("D23") "227.40"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +0.04%  "
Set-TextValue $ws.Range("D25") "2.40"
$ws.Range("E25").Value = "  -7.15%  "
$ws.Range("E26").Value = "  +0.56%  "
Set-TextValue $ws.Range("D27") "162.38"
$ws.Range("E27").Value = "  -0.25%  "
Set-TextValue $ws.Range("D28") "8.64"
$ws.Range("E28").Value = "  -4.44%  "
Set-TextValue $ws.Range("D29") "19.10"
$ws.Range("E29").Value = "  -3.43%  "
Set-TextValue $ws.Range("D30") "0.126"
$ws.Range("E30").Value = "  -3.96%  "
Set-TextValue $ws.Range("D32") "0.117"
$ws.Range("E32").Value = "  -2.20%  "
Set-TextValue $ws.Range("D33") "4.41"
$ws.Range("E33").Value = "  -5.57%  "
Set-TextValue $ws.Range("D34") "0.0608"
$ws.Range("E34").Value = "  -7.75%  "
Set-TextValue $ws.Range("D35") "4.22"
$ws.Range("E35").Value = "  -6.69%  "
Set-TextValue $ws.Range("D36") "2.29"
$ws.Range("E36").Value = "  -6.38%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -1.86%  "
Set-TextValue $ws.Range("D39") "3.24"
$ws.Range("E39").Value = "  -4.14%  "
Set-TextValue $ws.Range("D40") "5.30"
$ws.Range("E40").Value = "  +0.02%  "
Set-TextValue $ws.Range("D41") "3.04"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "1.425.95"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D43") "0.0203"
$ws.Range("E43").Value = "  -5.22%  "
Set-TextValue $ws.Range("D44") "1.12"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D45") "0.0892"
$ws.Range("E45").Value = "  -7.29%  "
Set-TextValue $ws.Range("D46") "87.80"
$ws.Range("E46").Value = "  -2.58%  "
Set-TextValue $ws.Range("D47") "1.00"
$ws.Range("E47").Value = "  -2.96%  "
Set-TextValue $ws.Range("D48") "14.95"
$ws.Range("E48").Value = "  -6.39%  "
Set-TextValue $ws.Range("D49") "2.89"
$ws.Range("E49").Value = "  +0.68%  "
Set-TextValue $ws.Range("D50") "6.72"
$ws.Range("E50").Value = "  -8.40%  "
$ws.Range("D51").Value = "2.176.00"
$ws.Range("E51").Value = "  -1.69%  "
